# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), columns AD:AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by copying
# the format from the last existing header cell (AC1).
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows 2-48: same team record for every row
for ($r = 2; $r -le 48; $r++) {
    $ws.Cells.Item($r, 30).Value = 106  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 56   # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF - Ties
}
